# Updates crypto price/volume data and reorders a few coin rows, per the
# "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''254.80'
$ws.Range("D3").Value = '''27.60'
$ws.Range("E3").Value = '''-7.95%'
$ws.Range("D4").Value = '''5.231'
$ws.Range("E4").Value = '''1.38%'
$ws.Range("D5").Value = '''0.05872'
$ws.Range("E5").Value = '''1.98%'
$ws.Range("D6").Value = '''6.720'
$ws.Range("E6").Value = '''0.81%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''3.219'
$ws.Range("E7").Value = '''-2.26%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.8636'
$ws.Range("E8").Value = '''1.51%'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = '''0.9660'
$ws.Range("E9").Value = '''12.52%'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '''0.0006090'
$ws.Range("E10").Value = '''1.81%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1413'
$ws.Range("E11").Value = '''1.96%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07164'
$ws.Range("E12").Value = '''1.25%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03186'
$ws.Range("E13").Value = '''-1.31%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09235'
$ws.Range("E14").Value = '''-1.40%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001544'
$ws.Range("E15").Value = '''0.73%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005823'
$ws.Range("E16").Value = '''-1.42%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.500'
$ws.Range("E17").Value = '''-1.33%'
$ws.Range("D18").Value = '''2.222'
$ws.Range("E18").Value = '''1.77%'
$ws.Range("D19").Value = '''0.3177'
$ws.Range("E19").Value = '''0.94%'
$ws.Range("D20").Value = '''0.03471'
$ws.Range("E20").Value = '''1.52%'
$ws.Range("D21").Value = '''0.1299'
$ws.Range("E21").Value = '''-1.46%'
$ws.Range("D22").Value = '''3.541'
$ws.Range("E22").Value = '''1.68%'
$ws.Range("D23").Value = '''0.04149'
$ws.Range("E23").Value = '''0.34%'
$ws.Range("E24").Value = '''-2.12%'
$ws.Range("E25").Value = '''-0.10%'
$ws.Range("D26").Value = '''0.004802'
$ws.Range("E26").Value = '''15.55%'
$ws.Range("E27").Value = '''0.00%'
$ws.Range("E28").Value = '''1.14%'
$ws.Range("D40").Value = '''0.03815'
$ws.Range("E40").Value = '''1.70%'
$ws.Range("D41").Value = '''0.005680'
$ws.Range("E41").Value = '''58.56%'
$ws.Range("D42").Value = '''0.1102'
$ws.Range("E42").Value = '''2.93%'
$ws.Range("E43").Value = '''-6.51%'
$ws.Range("D44").Value = '''0.01066'
$ws.Range("E44").Value = '''10.27%'
$ws.Range("D45").Value = '''0.00005242'
$ws.Range("E45").Value = '''-4.61%'
$ws.Range("E46").Value = '''-0.03%'
$ws.Range("D47").Value = '''0.09999'
$ws.Range("E47").Value = '''40.79%'
$ws.Range("D48").Value = '''0.002129'
$ws.Range("E48").Value = '''-14.00%'
$ws.Range("E49").Value = '''-0.03%'
$ws.Range("E50").Value = '''-0.03%'
